$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 10 ("Registration, ML Estimation, and Objective Function"): add a new
# content text box ("Content Placeholder 2") listing the topics for the
# baseline / conclusion discussion (variance of an estimator, Fisher
# information, CRLB, and its application to image registration).
# ---------------------------------------------------------------------------
$s10 = $p.Slides.Item(10)

# Position/size are expressed in points for AddTextbox (EMU / 12700):
#   off  x=457200  y=1219200  -> 36, 96
#   ext  cx=8229600 cy=4937760 -> 648, 388.8
$newShape = $s10.Shapes.AddTextbox(1, 36, 96, 648, 388.8)
$newShape.Name = "Content Placeholder 2"

$tr = $newShape.TextFrame.TextRange
$tr.Text = "Variance of an estimator`rFisher Information`rCramer-Rao lower bound (CRLB)`rQuantitative measure of estimator performance`rApplication of CRLB to image registration`r"

# "Quantitative measure of estimator performance" is a sub-bullet (outline
# level 2, i.e. OOXML lvl="1").
$tr.Paragraphs(4, 1).IndentLevel = 2

# ---------------------------------------------------------------------------
# Slide 3 ("Robotic Helicopter to Inspect Fukushima Reactors"): the body
# placeholder's shrink-text-to-fit scale was nudged from 77.5% to 85%.
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$bodyShape = $s3.Shapes.Item(2)
$bodyShape.TextFrame2.FontScale = 85000
